# Refresh the "cryptos" price table (coinranking.com snapshot).
# Mirrors the scheduled GitHub Actions job that re-scrapes the
# current price / 1h volume change for each coin.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    # The Price column stores every value as text in the source sheet
    # (even values that look like plain numbers, e.g. "0.994"), so pin
    # the cell to Text before writing, then restore the default/Normal
    # style so no stray number-format is left behind.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# --- Coin rows with updated Price (D) and Volume 1h (E) ---
$ws.Range("D2").Value = "66.350.62"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "3.350.76"
$ws.Range("E3").Value = "  +0.42%  "
Set-TextValue "D4" "0.994"
$ws.Range("E4").Value = "  -0.59%  "
Set-TextValue "D5" "586.38"
$ws.Range("E5").Value = "  +4.25%  "
Set-TextValue "D6" "186.00"
$ws.Range("E6").Value = "  -1.71%  "
Set-TextValue "D7" "0.999"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "3.345.77"
$ws.Range("E8").Value = "  +0.51%  "
Set-TextValue "D9" "0.576"
$ws.Range("E9").Value = "  -2.45%  "
Set-TextValue "D10" "0.182"
$ws.Range("E10").Value = "  -2.04%  "
Set-TextValue "D11" "0.584"
$ws.Range("E11").Value = "  -1.32%  "
Set-TextValue "D12" "47.14"
$ws.Range("E12").Value = "  -1.80%  "
Set-TextValue "D13" "0.0000269"
$ws.Range("E13").Value = "  -1.50%  "
Set-TextValue "D14" "668.91"
$ws.Range("E14").Value = "  +10.11%  "
$ws.Range("D15").Value = "3.865.16"
$ws.Range("E15").Value = "  -0.07%  "
Set-TextValue "D16" "8.53"
$ws.Range("E16").Value = "  -2.08%  "
$ws.Range("D17").Value = "66.386.60"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D20").Value = "3.336.28"
$ws.Range("E20").Value = "  -0.15%  "
Set-TextValue "D21" "11.13"
$ws.Range("E21").Value = "  -0.60%  "
Set-TextValue "D22" "0.901"
$ws.Range("E22").Value = "  -1.99%  "
Set-TextValue "D23" "17.83"
$ws.Range("E23").Value = "  -4.33%  "
Set-TextValue "D24" "102.17"
$ws.Range("E24").Value = "  +1.37%  "
Set-TextValue "D25" "5.05"
$ws.Range("E25").Value = "  -2.48%  "
Set-TextValue "D26" "3.99"
$ws.Range("E26").Value = "  -0.76%  "
Set-TextValue "D27" "2.79"
$ws.Range("E27").Value = "  +0.17%  "
Set-TextValue "D28" "9.46"
$ws.Range("E28").Value = "  -2.85%  "
Set-TextValue "D29" "32.29"
$ws.Range("E29").Value = "  +5.37%  "
Set-TextValue "D30" "8.53"
$ws.Range("E30").Value = "  -2.10%  "
Set-TextValue "D31" "6.87"
$ws.Range("E31").Value = "  +0.29%  "
Set-TextValue "D32" "611.83"
$ws.Range("E32").Value = "  +5.07%  "
Set-TextValue "D33" "3.92"
$ws.Range("E33").Value = "  +0.35%  "
Set-TextValue "D34" "11.13"
$ws.Range("E34").Value = "  -0.53%  "
$ws.Range("D35").Value = "3.855.86"
$ws.Range("E35").Value = "  +3.52%  "
Set-TextValue "D36" "0.106"
$ws.Range("E36").Value = "  -0.77%  "
Set-TextValue "D38" "56.28"
$ws.Range("E38").Value = "  -1.79%  "
Set-TextValue "D39" "0.129"
$ws.Range("E39").Value = "  -2.51%  "
Set-TextValue "D40" "2.68"
$ws.Range("E40").Value = "  -1.21%  "
$ws.Range("D41").Value = "0.0₃0704"
$ws.Range("E41").Value = "  -4.37%  "
Set-TextValue "D42" "33.09"
$ws.Range("E42").Value = "  -3.42%  "
Set-TextValue "D43" "3.21"
$ws.Range("E43").Value = "  -3.26%  "
Set-TextValue "D44" "3.40"
$ws.Range("E44").Value = "  +1.56%  "
Set-TextValue "D45" "0.338"
$ws.Range("E45").Value = "  -2.49%  "
Set-TextValue "D46" "0.0419"
$ws.Range("E46").Value = "  -2.28%  "
Set-TextValue "D47" "3.00"
$ws.Range("E47").Value = "  -13.95%  "

# --- Rows with only a Volume 1h (E) update ---
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("E48").Value = "  -1.87%  "
$ws.Range("E51").Value = "  +1.54%  "

# --- Rows whose Coin/Link/Price/Volume were re-ranked (swapped order) ---
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D18" "17.93"
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D19" "0.118"
$ws.Range("E19").Value = "  -0.68%  "
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue "D49" "2.56"
$ws.Range("E49").Value = "  -2.36%  "
$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D50" "1.00"
$ws.Range("E50").Value = "  +0.20%  "

Write-Host "Cryptos list updated"
